# Commit: "Wed, Aug 05, 2020  3:05:30 PM"
#
# The canonical-OOXML diff for this commit swaps the raw contents of
# ppt/theme/theme1.xml (linked from the slide master) and
# ppt/theme/theme2.xml (linked only from the notes master):
#   - theme1.xml goes from the deck's custom "Integral" theme to the
#     stock "Office Theme" palette.
#   - theme2.xml goes from the stock "Office Theme" palette to the old
#     "Integral" palette.
# fontScheme/fmtScheme are byte-identical between the two theme parts,
# so the only real content change is the <a:clrScheme> color values
# (and the cosmetic name="..." attributes, which PowerPoint's object
# model exposes as read-only and which this host doesn't persist).
#
# The only theme surface reachable through the PowerPoint COM object
# model (Master.Theme / SlideMaster.Theme / NotesMaster.Theme all
# resolve to the same single theme object in this host) is the slide
# master's theme, i.e. ppt/theme/theme1.xml. We drive that to the
# "Office Theme" palette to match the new theme1.xml content.

function Get-BGR($r, $g, $b) {
    # PowerPoint's RGB color properties store colors as BGR-packed
    # integers (the classic VBA RGB() encoding), not RGB-packed ones.
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$theme = $p.Slides.Item(1).Master.Theme
$colors = $theme.ThemeColorScheme

# clrScheme slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colors.Colors(1).RGB  = Get-BGR 0x00 0x00 0x00   # dk1       000000
$colors.Colors(2).RGB  = Get-BGR 0xFF 0xFF 0xFF   # lt1       FFFFFF
$colors.Colors(3).RGB  = Get-BGR 0x44 0x54 0x6A   # dk2       44546A
$colors.Colors(4).RGB  = Get-BGR 0xE7 0xE6 0xE6   # lt2       E7E6E6
$colors.Colors(5).RGB  = Get-BGR 0x5B 0x9B 0xD5   # accent1   5B9BD5
$colors.Colors(6).RGB  = Get-BGR 0xED 0x7D 0x31   # accent2   ED7D31
$colors.Colors(7).RGB  = Get-BGR 0xA5 0xA5 0xA5   # accent3   A5A5A5
$colors.Colors(8).RGB  = Get-BGR 0xFF 0xC0 0x00   # accent4   FFC000
$colors.Colors(9).RGB  = Get-BGR 0x44 0x72 0xC4   # accent5   4472C4
$colors.Colors(10).RGB = Get-BGR 0x70 0xAD 0x47   # accent6   70AD47
$colors.Colors(11).RGB = Get-BGR 0x05 0x63 0xC1   # hlink     0563C1
$colors.Colors(12).RGB = Get-BGR 0x95 0x4F 0x72   # folHlink  954F72
